$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pron")

# --- Rebuild the "pron" sheet's A1:I6 table -----------------------------
# Remove the old 5-column (A:E) layout and its merges, then lay out the
# new 9-column (A:I) layout that adds "与格" (dative) and "重读" (stressed)
# pronoun columns plus several new pronoun forms (me/te/le/leur/la).
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# Row 1 - "单数" (singular) / "复数" (plural) group headers
$ws.Range("A1").Value = ""
$ws.Range("B1").Value = "单数"
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = ""
$ws.Range("E1").Value = ""
$ws.Range("F1").Value = "复数"
$ws.Range("G1").Value = ""
$ws.Range("H1").Value = ""

# Row 2 - case headers: 主格/宾格/与格/重读 (repeated for each group)
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "主格"
$ws.Range("C2").Value = "宾格"
$ws.Range("D2").Value = "与格"
$ws.Range("E2").Value = "重读"
$ws.Range("F2").Value = "主格"
$ws.Range("G2").Value = "宾格"
$ws.Range("H2").Value = "与格"
$ws.Range("I2").Value = "重读"

# Row 3 - 第一人称 (1st person)
$ws.Range("A3").Value = "第一人称"
$ws.Range("B3").Value = "je"
$ws.Range("C3").Value = "me"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "moi"
$ws.Range("F3").Value = "nous"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

# Row 4 - 第二人称 (2nd person)
$ws.Range("A4").Value = "第二人称"
$ws.Range("B4").Value = "tu"
$ws.Range("C4").Value = "te"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "toi"
$ws.Range("F4").Value = "vous"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

# Row 5 - 第三人称 (3rd person), masculine
$ws.Range("A5").Value = "第三人称"
$ws.Range("B5").Value = "il"
$ws.Range("C5").Value = "le"
$ws.Range("D5").Value = "lui"
$ws.Range("E5").Value = "lui"
$ws.Range("F5").Value = "ils"
$ws.Range("G5").Value = "leur"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "eux"

# Row 6 - 第三人称 (3rd person), feminine
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "elle"
$ws.Range("C6").Value = "la"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "elle"
$ws.Range("F6").Value = "elles"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

# --- Alignment: every cell in the table is centered both ways ----------
$ws.Range("A1:H1").HorizontalAlignment = -4108
$ws.Range("A1:H1").VerticalAlignment = -4108
$ws.Range("A2:I6").HorizontalAlignment = -4108
$ws.Range("A2:I6").VerticalAlignment = -4108

# --- Merges --------------------------------------------------------------
$ws.Range("B1:E1").Merge()
$ws.Range("F1:H1").Merge()
$ws.Range("C3:D3").Merge()
$ws.Range("F3:I3").Merge()
$ws.Range("C4:D4").Merge()
$ws.Range("F4:I4").Merge()
$ws.Range("A1:A2").Merge()
$ws.Range("A5:A6").Merge()
$ws.Range("D5:D6").Merge()
$ws.Range("I5:I6").Merge()
$ws.Range("G5:H6").Merge()

# --- Column widths: C:E (was only C) match the new 3-column layout ------
$ws.Range("C1:E1").EntireColumn.ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Outline levels: keep row outline level, drop stray column outline --
$ws.Outline.ShowLevels(5, 0)

# --- Window / selection state --------------------------------------------
# Active tab moves from "prep-loc" back to "pron", with a new selection.
$ws.Activate()
$ws.Range("H16").Select()
